$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.652.66"
$ws.Range("E2").Value = "  +2.89%  "

$ws.Range("D3").Value = "1.697.31"
$ws.Range("E3").Value = "  +2.40%  "

$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  -0.90%  "

$ws.Range("D5").Value = "'314.01"
$ws.Range("E5").Value = "  +1.69%  "

$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = "  -0.78%  "

$ws.Range("D7").Value = "'0.3956"
$ws.Range("E7").Value = "  +1.67%  "

$ws.Range("D8").Value = "'0.4047"
$ws.Range("E8").Value = "  +3.29%  "

$ws.Range("D9").Value = "'57.48"
$ws.Range("E9").Value = "  +20.55%  "

$ws.Range("D10").Value = "'1.527"
$ws.Range("E10").Value = "  +10.65%  "

$ws.Range("D11").Value = "'1.000"
$ws.Range("E11").Value = "  -0.87%  "

$ws.Range("D12").Value = "'0.08783"
$ws.Range("E12").Value = "  +2.31%  "

$ws.Range("D13").Value = "'7.309"
$ws.Range("E13").Value = "  +14.01%  "

$ws.Range("D14").Value = "'23.19"
$ws.Range("E14").Value = "  +3.12%  "

$ws.Range("E15").Value = "  +2.58%  "

$ws.Range("D16").Value = "'7.628"
$ws.Range("E16").Value = "  +7.35%  "

$ws.Range("D17").Value = "1.695.35"
$ws.Range("E17").Value = "  +2.09%  "

$ws.Range("D18").Value = "'100.48"
$ws.Range("E18").Value = "  +0.47%  "

$ws.Range("D19").Value = "'0.07060"
$ws.Range("E19").Value = "  +4.60%  "

$ws.Range("D20").Value = "'19.51"
$ws.Range("E20").Value = "  +3.70%  "

$ws.Range("D21").Value = "'6.734"
$ws.Range("E21").Value = "  +2.56%  "

$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "  -0.73%  "

$ws.Range("D23").Value = "'14.19"
$ws.Range("E23").Value = "  +4.66%  "

$ws.Range("D24").Value = "24.592.40"
$ws.Range("E24").Value = "  +2.61%  "

$ws.Range("D25").Value = "'3.005"
$ws.Range("E25").Value = "  +13.24%  "

$ws.Range("D26").Value = "'2.311"
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").Value = "'22.44"
$ws.Range("E27").Value = "  +3.76%  "

$ws.Range("D28").Value = "'159.68"
$ws.Range("E28").Value = "  +1.62%  "

$ws.Range("D29").Value = "'5.190"
$ws.Range("E29").Value = "  +1.71%  "

$ws.Range("D30").Value = "'133.57"
$ws.Range("E30").Value = "  +3.34%  "

$ws.Range("D31").Value = "'7.571"
$ws.Range("E31").Value = "  +37.02%  "

$ws.Range("D32").Value = "1.881.17"
$ws.Range("E32").Value = "  +1.98%  "

$ws.Range("D33").Value = "'1.089"
$ws.Range("E33").Value = "  -1.32%  "

$ws.Range("D34").Value = "'7.349"
$ws.Range("E34").Value = "  +21.38%  "

$ws.Range("D35").Value = "'0.08562"
$ws.Range("E35").Value = "  +1.28%  "

$ws.Range("D36").Value = "'1.961"
$ws.Range("E36").Value = "  +11.07%  "

$ws.Range("D37").Value = "'11.02"
$ws.Range("E37").Value = "  +7.09%  "

$ws.Range("D38").Value = "'0.2719"
$ws.Range("E38").Value = "  +4.58%  "

$ws.Range("D39").Value = "'14.75"
$ws.Range("E39").Value = "  -0.52%  "

$ws.Range("D40").Value = "'0.02778"
$ws.Range("E40").Value = "  +12.25%  "

$ws.Range("D41").Value = "'0.09056"
$ws.Range("E41").Value = "  +3.49%  "

$ws.Range("D42").Value = "'1.475"
$ws.Range("E42").Value = "  +3.67%  "

$ws.Range("D43").Value = "'0.7642"
$ws.Range("E43").Value = "  +4.97%  "

$ws.Range("D44").Value = "'0.7187"
$ws.Range("E44").Value = "  +4.47%  "

$ws.Range("D45").Value = "'15.40"
$ws.Range("E45").Value = "  +5.06%  "

$ws.Range("D46").Value = "'2.461"
$ws.Range("E46").Value = "  +5.17%  "

$ws.Range("D47").Value = "'4.166"
$ws.Range("E47").Value = "  +2.58%  "

$ws.Range("D48").Value = "'0.9988"
$ws.Range("E48").Value = "  -0.84%  "

$ws.Range("D49").Value = "'1.329"
$ws.Range("E49").Value = "  +19.39%  "

$ws.Range("D50").Value = "'140.46"
$ws.Range("E50").Value = "  +1.68%  "

$ws.Range("D51").Value = "'0.00000000380"
$ws.Range("E51").Value = "  +2.25%  "

Write-Host "Updated cryptos list"
